$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$t1 = $s1.Shapes.Item(1).TextFrame.TextRange
$t1.Text = " "
$t1.Text = "Example numbering MWE"

$s2 = $p.Slides.Item(2)
$t2 = $s2.Shapes.Item(1).TextFrame.TextRange
$t2.Text = " "
$t2.Text = "A second slide"
